$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Easter Sunday" caption text (P38) ---
$ws.Range("P38").Value = "EASTER SUNDAY - THE DAY TRUMP WANTS CORONAVIRUS TO GO AWAY"

# --- Update the view/window scroll position (best effort; cosmetic) ---
$win = $wb.Windows.Item(1)
$win.Left = 28720
$win.Top = 1860
$win.ScrollRow = 15
$win.ScrollColumn = 1

# --- Replace I23's forecast formula with the actual reported case count for 3/28 ---
# First copy I22's "actual" formatting (fill/number format) onto I23 so the cell
# switches from the "forecast" style to the "actual" style, matching the other
# hard-coded actual cells (I20:I22).
$ws.Range("I22").Copy()
$ws.Range("I23").PasteSpecial(-4122)
$ws.Range("I23").Value = 123578

# I24 now becomes the first cell of the forecast chain that follows the new
# actual value, recompute it from I23 directly.
$ws.Range("I24").Formula = "=I23*(1+AVERAGE(M22:M23))"
